$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")
foreach ($w in 40,40.4,40.5,40.6,41,44.17,44.5,44.6) {
  $ws3.Columns.Item(2).ColumnWidth = $w
  Write-Host ("set=" + $w + " got=" + $ws3.Columns.Item(2).ColumnWidth)
}
